$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column C, shifting Mem Capacitance etc. one column right
$ws.Columns.Item(3).EntireColumn.Insert()

# Populate the new column C headers
$ws.Range("C1").Value = "Initial Voltage"
$ws.Range("C2").Value = "U0"
$ws.Range("C3").Value = "[V]"

# Populate the new column C data values (Initial Voltage = 0 for each CPG)
$ws.Range("C4").Value = 0
$ws.Range("C5").Value = 0
$ws.Range("C6").Value = 0
$ws.Range("C7").Value = 0

# Match column width to neighboring column (same as old column C's width).
# (ColumnWidth is quantized in sixths of a character by this engine, so 16
# is the closest achievable value to the original 16.85546875 stored width.)
$ws.Columns.Item(3).ColumnWidth = 16

# Update selection as shown in the diff
$ws.Range("C4").Select() | Out-Null
